# 06_LibFormula.xlsx - "Indicatori BIB complessita 2 Roberto"
#
# Sheet "Library_Formula": column C (rows 28..80) holds a list of
# INDICATOR_<n> tokens sorted ascending by <n>. This change adds 19 new
# indicators to that list and re-sorts the whole run, which now spans
# rows 28..99. Columns A/B/E are repeated constants on every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")

# Extend the sheet with 19 new, correctly-styled rows (inherits the
# existing per-column styles instead of defaulting to column styles).
$ws.Range("A81:G99").Insert()

# Final ascending-sorted indicator numbers for rows 28..99 (53 existing +
# 19 new: 40, 57, 71, 77, 81, 85, 89, 93, 98, 102, 106, 110, 114, 119,
# 125, 129, 133, 181, 201).
$indicatorNumbers = @(40,48,51,56,57,60,64,66,69,70,71,74,75,76,77,78,79,80,81,82,83,84,85,86,87,88,89,90,91,92,93,94,95,96,97,98,99,100,102,103,104,105,106,107,108,110,111,112,113,114,115,116,118,119,120,121,124,125,126,127,128,129,130,131,132,133,173,181,182,186,201,205)

$startRow = 28
for ($i = 0; $i -lt $indicatorNumbers.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("A$row").Value = "CREATE/MODIFY"
    $ws.Range("B$row").Value = "LIB_EWS_BE"
    $ws.Range("C$row").Value = "INDICATOR_" + $indicatorNumbers[$i]
    $ws.Range("E$row").Value = "String"
}

# Reflect the scrolled/selected state captured in the saved workbook.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 83
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E80:E99").Select() | Out-Null
